# Applies the cryptos-list price/volume refresh described in the commit diff.
# Most D/E cells are plain text updates; a subset of D-column price cells hold
# strings that LOOK like numbers (e.g. "1.000", "0.08800") and must be forced to
# stay text (matching the original inlineStr cells) instead of being coerced to
# doubles by the COM Value setter - hence the NumberFormat/Style dance below.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.759.01'
$ws.Range('E2').Value = '  -1.47%  '
$ws.Range('D3').Value = '1.798.01'
$ws.Range('E3').Value = '  -1.33%  '
$ws.Range('E4').Value = '  -0.34%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.12'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.94%  '
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4460'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +5.65%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3667'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07322'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.30%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8554'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.58'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Value = '1.919.20'
$ws.Range('E12').Value = '  +4.96%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.605'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '92.06'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.75%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07072'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.07%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.284'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.23%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008688'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.70%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.000'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.33%  '
$ws.Range('E20').Value = '  -1.18%  '
$ws.Range('D21').Value = '26.803.89'
$ws.Range('E21').Value = '  -1.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.148'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.69%  '
$ws.Range('E23').Value = '  -0.71%  '
$ws.Range('E24').Value = '  +0.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.97'
$ws.Range('D25').Style = 'Normal'
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.182'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.61%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.45'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.50%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.191'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E29').Value = '  +0.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08800'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.7441'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.157'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.938'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.449'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.001'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.084'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01958'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05172'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5292'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.846'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.68%  '
$ws.Range('E41').Value = '  -3.40%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1682'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.73%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5100'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.75%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.391'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.54'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.962'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.33%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '105.53'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.78%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.001'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.33%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.659'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06304'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9145'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.35%  '
